$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 182, shifting existing
# rows 182-190 down to rows 184-192 (preserving their values/styles).
$ws.Rows("182:183").Insert()

# --- New row 182 ---
$ws.Range("A182").Value = 5
$ws.Range("B182").Value = "Macroferia Regional de Talca"
$ws.Range("C182").Value = "Maule"
$ws.Range("D182").Value = 44610
$ws.Range("E182").Value = 7
$ws.Range("F182").Value = 100112021
$ws.Range("G182").Value = "Ají"
$ws.Range("H182").Value = "Cacho cabra verde"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 150
$ws.Range("K182").Value = 15000
$ws.Range("L182").Value = 15000
$ws.Range("M182").Value = 15000
$ws.Range("N182").Value = "$/saco 25 kilos"
$ws.Range("O182").Value = "Región del Maule"
$ws.Range("P182").Value = 600
$ws.Range("Q182").Value = 25
$ws.Range("R182").Value = "Hortaliza"

# --- New row 183 ---
$ws.Range("A183").Value = 5
$ws.Range("B183").Value = "Macroferia Regional de Talca"
$ws.Range("C183").Value = "Maule"
$ws.Range("D183").Value = 44610
$ws.Range("E183").Value = 7
$ws.Range("F183").Value = 100112021
$ws.Range("G183").Value = "Ají"
$ws.Range("H183").Value = "Cristal"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 150
$ws.Range("K183").Value = 15000
$ws.Range("L183").Value = 15000
$ws.Range("M183").Value = 15000
$ws.Range("N183").Value = "$/saco 25 kilos"
$ws.Range("O183").Value = "Región del Maule"
$ws.Range("P183").Value = 600
$ws.Range("Q183").Value = 25
$ws.Range("R183").Value = "Hortaliza"
